# Update "number of visits / views" (column F) values for several rows
# across three worksheets, per the latest site data regeneration.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheetId 1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 198
$ws1.Range("F8").Value  = 13133
$ws1.Range("F9").Value  = 2270
$ws1.Range("F12").Value = 53926
$ws1.Range("F13").Value = 1294
$ws1.Range("F14").Value = 310
$ws1.Range("F16").Value = 854
$ws1.Range("F20").Value = 852
$ws1.Range("F22").Value = 1247
$ws1.Range("F28").Value = 1190
$ws1.Range("F36").Value = 43
$ws1.Range("F37").Value = 4698
$ws1.Range("F38").Value = 32
$ws1.Range("F40").Value = 8692
$ws1.Range("F46").Value = 97

# --- Sheet "本地生活" (sheetId 3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 30

# --- Sheet "全部类型" (sheetId 4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 198
$ws4.Range("F8").Value  = 13133
$ws4.Range("F9").Value  = 13133
$ws4.Range("F10").Value = 2270
$ws4.Range("F11").Value = 1294
$ws4.Range("F12").Value = 310
$ws4.Range("F13").Value = 854
$ws4.Range("F17").Value = 852
$ws4.Range("F19").Value = 1247
$ws4.Range("F20").Value = 30
$ws4.Range("F26").Value = 1190
$ws4.Range("F33").Value = 43
$ws4.Range("F34").Value = 4698
$ws4.Range("F36").Value = 8692
$ws4.Range("F42").Value = 97

$wb.Save()
